# add second method to price scenario generation - moving average to de-trend data.
# This updates the pre-computed price-scenario values on the "grain", "meat", "wool"
# and "prob" sheets to reflect the new (moving-average de-trended) calculation.

$wb = $excel.ActiveWorkbook

$columnsBtoK = @("B","C","D","E","F","G","H","I","J","K")

# Sheets "grain", "meat" and "wool" each hold a 4-row x 10-column (B:K) block of
# scenario multipliers; every cell in a given row shares the same value.

$wsGrain = $wb.Worksheets.Item("grain")
$rowValuesGrain = @{
    2 = 0.9001469258456177
    3 = 0.8759941262449619
    4 = 1.124005873755038
    5 = 1.099853074154382
}
foreach ($rowNum in $rowValuesGrain.Keys) {
    $value = $rowValuesGrain[$rowNum]
    foreach ($col in $columnsBtoK) {
        $wsGrain.Range("$col$rowNum").Value = $value
    }
}

$wsMeat = $wb.Worksheets.Item("meat")
$rowValuesMeat = @{
    2 = 0.8167959426209745
    3 = 1.227518074963115
    4 = 0.7724819250368852
    5 = 1.183204057379025
}
foreach ($rowNum in $rowValuesMeat.Keys) {
    $value = $rowValuesMeat[$rowNum]
    foreach ($col in $columnsBtoK) {
        $wsMeat.Range("$col$rowNum").Value = $value
    }
}

$wsWool = $wb.Worksheets.Item("wool")
$rowValuesWool = @{
    2 = 0.8167959426209745
    3 = 1.227518074963115
    4 = 0.7724819250368852
    5 = 1.183204057379025
}
foreach ($rowNum in $rowValuesWool.Keys) {
    $value = $rowValuesWool[$rowNum]
    foreach ($col in $columnsBtoK) {
        $wsWool.Range("$col$rowNum").Value = $value
    }
}

# Sheet "prob" holds a single column (B) of probabilities, one per row.
$wsProb = $wb.Worksheets.Item("prob")
$rowValuesProb = @{
    2 = 0.2023576762861144
    3 = 0.2976423237138854
    4 = 0.2976423237138855
    5 = 0.2023576762861146
}
foreach ($rowNum in $rowValuesProb.Keys) {
    $value = $rowValuesProb[$rowNum]
    $wsProb.Range("B$rowNum").Value = $value
}
